$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 10004
$ws.Range("I10").Value = 10004
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 10004
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -9711
$ws.Range("N10").ClearContents()

$ws.Range("H11").Value = 177.33333
$ws.Range("I11").Value = 177.33333
$ws.Range("K11").Value = 177.33333
$ws.Range("M11").Value = -37.33332999999999

$ws.Range("H15").Value = 1341.4878
$ws.Range("I15").Value = 1341.4878
$ws.Range("K15").Value = 4024.463400000001
$ws.Range("M15").Value = -3855.463400000001

$ws.Range("H17").Value = 1952.15
$ws.Range("J17").Value = 2551.4167
$ws.Range("L17").Value = 7654.250100000001
$ws.Range("N17").Value = -7990.250100000001

$ws.Range("H62").Value = 8061.6875
$ws.Range("I62").Value = 5122
$ws.Range("K62").Value = 5122
$ws.Range("M62").Value = -4498

$ws.Range("H65").Value = 8061.6875
$ws.Range("I65").Value = 5122
$ws.Range("K65").Value = 25610
$ws.Range("M65").Value = -22490

$ws.Range("H69").Value = 250004370
$ws.Range("J69").Value = 250004370
$ws.Range("L69").Value = 750013110
$ws.Range("N69").Value = -750014858

$ws.Range("H72").Value = 250004370
$ws.Range("J72").Value = 250004370
$ws.Range("L72").Value = 2250039330
$ws.Range("N72").Value = -2250048066

$ws.Range("H86").Value = 1722.2222
$ws.Range("I86").Value = 1880
$ws.Range("K86").Value = 1880
$ws.Range("M86").Value = -757

$ws.Range("H89").Value = 1722.2222
$ws.Range("I89").Value = 1880
$ws.Range("K89").Value = 9400
$ws.Range("M89").Value = -3784

$ws.Range("H92").Value = 1169.4762
$ws.Range("I92").Value = 974.1177
$ws.Range("K92").Value = 974.1177
$ws.Range("M92").Value = 273.8823

$ws.Range("H98").Value = 2644.1956
$ws.Range("I98").Value = 1678.8536
$ws.Range("J98").Value = 10560
$ws.Range("K98").Value = 1678.8536
$ws.Range("L98").Value = 10560
$ws.Range("M98").Value = -180.8535999999999
$ws.Range("N98").Value = -13556

$ws.Range("H100").Value = 7084.8
$ws.Range("I100").Value = 6632.6665
$ws.Range("J100").Value = 7278.5713
$ws.Range("K100").Value = 6632.6665
$ws.Range("L100").Value = 7278.5713
$ws.Range("M100").Value = -6091.6665
$ws.Range("N100").Value = -8360.5713

$ws.Range("H107").Value = 1610.1666
$ws.Range("J107").Value = 1146
$ws.Range("L107").Value = 1146
$ws.Range("N107").Value = -4986

$ws.Range("H111").Value = 2498
$ws.Range("I111").Value = 2498
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 7494
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -4427
$ws.Range("N111").ClearContents()

$ws.Range("H113").Value = 3828.8845
$ws.Range("I113").Value = 3819.652
$ws.Range("K113").Value = 3819.652
$ws.Range("M113").Value = -565.652

$ws.Range("H115").Value = 4653
$ws.Range("I115").Value = 3481.25
$ws.Range("J115").Value = 5824.75
$ws.Range("K115").Value = 10443.75
$ws.Range("L115").Value = 17474.25
$ws.Range("M115").Value = -8876.75
$ws.Range("N115").Value = -20608.25

$ws.Range("H116").Value = 5791.1333
$ws.Range("I116").Value = 3232.5557
$ws.Range("J116").Value = 9629
$ws.Range("K116").Value = 3232.5557
$ws.Range("L116").Value = 9629
$ws.Range("M116").Value = 209.4443000000001
$ws.Range("N116").Value = -16513

$ws.Range("H122").Value = 2644.1956
$ws.Range("I122").Value = 1678.8536
$ws.Range("J122").Value = 10560
$ws.Range("K122").Value = 5036.560799999999
$ws.Range("L122").Value = 31680
$ws.Range("M122").Value = -2586.560799999999
$ws.Range("N122").Value = -36580

$ws.Range("H123").Value = 67470.11
$ws.Range("J123").Value = 67470.11
$ws.Range("L123").Value = 67470.11
$ws.Range("N123").Value = -77270.11

$ws.Range("H131").Value = 7580894
$ws.Range("I131").Value = 12822067
$ws.Range("K131").Value = 38466201
$ws.Range("M131").Value = -38461161

$ws.Range("H132").Value = 2129832.5
$ws.Range("I132").Value = 2093.025
$ws.Range("J132").Value = 14288344
$ws.Range("K132").Value = 6279.075000000001
$ws.Range("L132").Value = 42865032
$ws.Range("M132").Value = -3749.075000000001
$ws.Range("N132").Value = -42870092

$ws.Range("H137").Value = 725783.4
$ws.Range("I137").Value = 966.1905
$ws.Range("J137").Value = 2417023.5
$ws.Range("K137").Value = 2898.5715
$ws.Range("L137").Value = 7251070.5
$ws.Range("M137").Value = -348.5715
$ws.Range("N137").Value = -7256170.5

$ws.Range("H138").Value = 1713196.8
$ws.Range("I138").Value = 1827.6
$ws.Range("J138").Value = 2226607.5
$ws.Range("K138").Value = 5482.799999999999
$ws.Range("L138").Value = 6679822.5
$ws.Range("M138").Value = -342.7999999999993
$ws.Range("N138").Value = -6690102.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1658.5333
$ws.Range("I2").Value = 443.45456
$ws.Range("K2").Value = 443.45456
$ws.Range("M2").Value = -330.45456

$ws.Range("H32").Value = 135070.84
$ws.Range("I32").Value = 135070.84
$ws.Range("K32").Value = 135070.84
$ws.Range("M32").Value = -134783.84

$ws.Range("H61").Value = 573250.5
$ws.Range("I61").Value = 5132.316
$ws.Range("J61").Value = 4171332.2
$ws.Range("K61").Value = 5132.316
$ws.Range("L61").Value = 4171332.2
$ws.Range("M61").Value = -4920.316
$ws.Range("N61").Value = -4171756.2

$ws.Range("H88").Value = 2780.5908
$ws.Range("J88").Value = 2584.6
$ws.Range("L88").Value = 2584.6
$ws.Range("N88").Value = -3396.6

$ws.Range("H91").Value = 2780.5908
$ws.Range("J91").Value = 2584.6
$ws.Range("L91").Value = 2584.6
$ws.Range("N91").Value = -5392.6

$ws.Range("H97").Value = 467.5
$ws.Range("I97").Value = 467.5
$ws.Range("K97").Value = 467.5
$ws.Range("M97").Value = 28.5

$ws.Range("H102").Value = 3236.7715
$ws.Range("I102").Value = 1100.3182
$ws.Range("K102").Value = 1100.3182
$ws.Range("M102").Value = 521.6818000000001

$ws.Range("H110").Value = 5230.1665
$ws.Range("I110").Value = 5230.1665
$ws.Range("K110").Value = 5230.1665
$ws.Range("M110").Value = -3185.1665

$ws.Range("H112").Value = 46166.668
$ws.Range("J112").Value = 46166.668
$ws.Range("L112").Value = 46166.668
$ws.Range("N112").Value = -49120.668

$ws.Range("H116").Value = 1658.5333
$ws.Range("I116").Value = 443.45456
$ws.Range("K116").Value = 443.45456
$ws.Range("M116").Value = 1850.54544

$ws.Range("H122").Value = 5566907.5
$ws.Range("I122").Value = 7415110
$ws.Range("J122").Value = 22299.6
$ws.Range("K122").Value = 22245330
$ws.Range("L122").Value = 66898.79999999999
$ws.Range("M122").Value = -22242880
$ws.Range("N122").Value = -71798.79999999999

$ws.Range("H132").Value = 5161.5
$ws.Range("I132").Value = 5161.5
$ws.Range("K132").Value = 15484.5
$ws.Range("M132").Value = -12954.5

$ws.Range("H136").Value = 573250.5
$ws.Range("I136").Value = 5132.316
$ws.Range("J136").Value = 4171332.2
$ws.Range("K136").Value = 15396.948
$ws.Range("L136").Value = 12513996.6
$ws.Range("M136").Value = -12846.948
$ws.Range("N136").Value = -12519096.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1658.5333
$ws.Range("I3").Value = 443.45456
$ws.Range("K3").Value = 443.45456
$ws.Range("M3").Value = -329.45456

$ws.Range("H16").Value = 418.18182
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H17").Value = 13250
$ws.Range("I17").Value = 10000
$ws.Range("J17").Value = 14333.333
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 14333.333
$ws.Range("M17").Value = -9828
$ws.Range("N17").Value = -14677.333

$ws.Range("H35").Value = 69997.5
$ws.Range("J35").Value = 69997.5
$ws.Range("L35").Value = 69997.5
$ws.Range("N35").Value = -70617.5

$ws.Range("H92").Value = 37436.363
$ws.Range("J92").Value = 37436.363
$ws.Range("L92").Value = 37436.363
$ws.Range("N92").Value = -42428.363

$ws.Range("H99").Value = 3135.2917
$ws.Range("J99").Value = 4271.375
$ws.Range("L99").Value = 4271.375
$ws.Range("N99").Value = -7267.375

$ws.Range("H134").Value = 405012.9
$ws.Range("I134").Value = 1793.1273
$ws.Range("K134").Value = 5379.3819
$ws.Range("M134").Value = -2844.3819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 4604.385
$ws.Range("I5").Value = 4714.5454
$ws.Range("J5").Value = 3998.5
$ws.Range("K5").Value = 4714.5454
$ws.Range("L5").Value = 3998.5
$ws.Range("M5").Value = -4602.5454
$ws.Range("N5").Value = -4222.5

$ws.Range("H12").Value = 1037.25
$ws.Range("I12").Value = 1037.25
$ws.Range("K12").Value = 1037.25
$ws.Range("M12").Value = -867.25

$ws.Range("H16").Value = 2692.8667
$ws.Range("I16").Value = 1883.1666
$ws.Range("J16").Value = 3232.6667
$ws.Range("K16").Value = 1883.1666
$ws.Range("L16").Value = 3232.6667
$ws.Range("M16").Value = -1596.1666
$ws.Range("N16").Value = -3806.6667

$ws.Range("H25").Value = 7755
$ws.Range("I25").Value = 1010
$ws.Range("J25").Value = 14500
$ws.Range("K25").Value = 1010
$ws.Range("L25").Value = 14500
$ws.Range("M25").Value = -836
$ws.Range("N25").Value = -14848

$ws.Range("H31").Value = 3331.077
$ws.Range("I31").Value = 1221.7
$ws.Range("J31").Value = 4649.4375
$ws.Range("K31").Value = 1221.7
$ws.Range("L31").Value = 4649.4375
$ws.Range("M31").Value = -926.7
$ws.Range("N31").Value = -5239.4375

$ws.Range("H34").Value = 3331.077
$ws.Range("I34").Value = 1221.7
$ws.Range("J34").Value = 4649.4375
$ws.Range("K34").Value = 1221.7
$ws.Range("L34").Value = 4649.4375
$ws.Range("M34").Value = -1019.7
$ws.Range("N34").Value = -5053.4375

$ws.Range("H37").Value = 1
$ws.Range("J37").Value = 1
$ws.Range("L37").Value = 1
$ws.Range("N37").Value = -215

$ws.Range("H58").Value = 1994.9584
$ws.Range("I58").Value = 2033.9131
$ws.Range("J58").Value = 1099
$ws.Range("K58").Value = 2033.9131
$ws.Range("L58").Value = 1099
$ws.Range("M58").Value = -1830.9131
$ws.Range("N58").Value = -1505

$ws.Range("H63").Value = 39999.9
$ws.Range("J63").Value = 39999.9
$ws.Range("L63").Value = 39999.9
$ws.Range("N63").Value = -41371.9

$ws.Range("H66").Value = 39999.9
$ws.Range("J66").Value = 39999.9
$ws.Range("L66").Value = 119999.7
$ws.Range("N66").Value = -126863.7

$ws.Range("H99").Value = 2282.1904
$ws.Range("I99").Value = 2206.5
$ws.Range("K99").Value = 2206.5
$ws.Range("M99").Value = -708.5

$ws.Range("H113").Value = 2692.8667
$ws.Range("I113").Value = 1883.1666
$ws.Range("J113").Value = 3232.6667
$ws.Range("K113").Value = 1883.1666
$ws.Range("L113").Value = 3232.6667
$ws.Range("M113").Value = 286.8334
$ws.Range("N113").Value = -7572.6667

$ws.Range("H126").Value = 2282.1904
$ws.Range("I126").Value = 2206.5
$ws.Range("K126").Value = 6619.5
$ws.Range("M126").Value = -4149.5

$ws.Range("H132").Value = 1696.75
$ws.Range("I132").Value = 952.625
$ws.Range("J132").Value = 7649.75
$ws.Range("K132").Value = 2857.875
$ws.Range("L132").Value = 22949.25
$ws.Range("M132").Value = -327.875
$ws.Range("N132").Value = -28009.25

$ws.Range("H134").Value = 1457.9231
$ws.Range("I134").Value = 1255.9714
$ws.Range("K134").Value = 3767.9142
$ws.Range("M134").Value = -1232.9142

$ws.Range("H136").Value = 1994.9584
$ws.Range("I136").Value = 2033.9131
$ws.Range("J136").Value = 1099
$ws.Range("K136").Value = 6101.7393
$ws.Range("L136").Value = 3297
$ws.Range("M136").Value = -3551.7393
$ws.Range("N136").Value = -8397

$ws.Range("H138").Value = 97390
$ws.Range("J138").Value = 97390
$ws.Range("L138").Value = 97390
$ws.Range("N138").Value = -107670

$ws.Range("H139").Value = 49999.5
$ws.Range("J139").Value = 49999.5
$ws.Range("L139").Value = 49999.5
$ws.Range("N139").Value = -60279.5

$ws.Range("H141").Value = 240609.89
$ws.Range("J141").Value = 240609.89
$ws.Range("L141").Value = 240609.89
$ws.Range("N141").Value = -250969.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1051.6111
$ws.Range("I5").Value = 749.06665
$ws.Range("K5").Value = 2247.19995
$ws.Range("M5").Value = -2135.19995

$ws.Range("H8").Value = 8810.134
$ws.Range("I8").Value = 8810.134
$ws.Range("K8").Value = 26430.402
$ws.Range("M8").Value = -26291.402

$ws.Range("H23").Value = 2096.3076
$ws.Range("J23").Value = 652.7273
$ws.Range("L23").Value = 1958.1819
$ws.Range("N23").Value = -2428.1819

$ws.Range("H33").Value = 309.54544
$ws.Range("I33").Value = 161
$ws.Range("J33").Value = 433.33334
$ws.Range("K33").Value = 966
$ws.Range("L33").Value = 2600.00004
$ws.Range("M33").Value = -683
$ws.Range("N33").Value = -3166.00004

$ws.Range("H34").Value = 3779.743
$ws.Range("J34").Value = 3779.743
$ws.Range("L34").Value = 11339.229
$ws.Range("N34").Value = -11507.229

$ws.Range("H82").Value = 66800
$ws.Range("J82").Value = 66800
$ws.Range("L82").Value = 200400
$ws.Range("N82").Value = -201212

$ws.Range("H85").Value = 66800
$ws.Range("J85").Value = 66800
$ws.Range("L85").Value = 200400
$ws.Range("N85").Value = -203208

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 2552
$ws.Range("J107").Value = 2043.9333
$ws.Range("L107").Value = 6131.7999
$ws.Range("N107").Value = -9971.7999

$ws.Range("H110").Value = 5010.6
$ws.Range("I110").Value = 5010.6
$ws.Range("K110").Value = 15031.8
$ws.Range("M110").Value = -10941.8

$ws.Range("H121").Value = 1553.8334
$ws.Range("I121").Value = 1418.3334
$ws.Range("K121").Value = 4255.0002
$ws.Range("M121").Value = -2945.0002

$ws.Range("H132").Value = 1598.9615
$ws.Range("I132").Value = 923.41174
$ws.Range("J132").Value = 2875
$ws.Range("K132").Value = 8310.70566
$ws.Range("L132").Value = 25875
$ws.Range("M132").Value = -5780.70566
$ws.Range("N132").Value = -30935

$ws.Range("H135").Value = 1051.6111
$ws.Range("I135").Value = 749.06665
$ws.Range("K135").Value = 6741.59985
$ws.Range("M135").Value = -4206.59985

$ws.Range("H136").Value = 3600.6428
$ws.Range("I136").Value = 2117.4167
$ws.Range("K136").Value = 6352.250100000001
$ws.Range("M136").Value = -1252.250100000001

$ws.Range("H137").Value = 16021.083
$ws.Range("I137").Value = 1452
$ws.Range("K137").Value = 4356
$ws.Range("M137").Value = 744

$ws.Range("H140").Value = 4830.7144
$ws.Range("I140").Value = 2420.8572
$ws.Range("J140").Value = 9650.429
$ws.Range("K140").Value = 7262.571599999999
$ws.Range("L140").Value = 28951.287
$ws.Range("M140").Value = -2082.571599999999
$ws.Range("N140").Value = -39311.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 102570.4
$ws.Range("I3").Value = 2931.4285
$ws.Range("J3").Value = 335061.34
$ws.Range("K3").Value = 2931.4285
$ws.Range("L3").Value = 335061.34
$ws.Range("M3").Value = -2815.4285
$ws.Range("N3").Value = -335293.34

$ws.Range("H11").Value = 1076846.5
$ws.Range("I11").Value = 1961750
$ws.Range("J11").Value = 191943
$ws.Range("K11").Value = 1961750
$ws.Range("L11").Value = 191943
$ws.Range("M11").Value = -1961611
$ws.Range("N11").Value = -192221

$ws.Range("H26").Value = 79999
$ws.Range("J26").Value = 79999
$ws.Range("L26").Value = 79999
$ws.Range("N26").Value = -80559

$ws.Range("H50").Value = 79999
$ws.Range("J50").Value = 79999
$ws.Range("L50").Value = 79999
$ws.Range("N50").Value = -80995

$ws.Range("H80").Value = 111114424
$ws.Range("I80").Value = 142859980
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 142859980
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -142858982
$ws.Range("N80").Value = -6996

$ws.Range("H83").Value = 111114424
$ws.Range("I83").Value = 142859980
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 714299900
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -714294908
$ws.Range("N83").Value = -34984

$ws.Range("H102").Value = 1018.9474
$ws.Range("I102").Value = 631.1111
$ws.Range("J102").Value = 8000
$ws.Range("K102").Value = 631.1111
$ws.Range("L102").Value = 8000
$ws.Range("M102").Value = 990.8889
$ws.Range("N102").Value = -11244

$ws.Range("H105").Value = 32668.834
$ws.Range("J105").Value = 32668.834
$ws.Range("L105").Value = 32668.834
$ws.Range("N105").Value = -39656.834

$ws.Range("H132").Value = 66668870
$ws.Range("J132").Value = 2857
$ws.Range("L132").Value = 8571
$ws.Range("N132").Value = -13631

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3016.52
$ws.Range("I7").Value = 2703.9524
$ws.Range("K7").Value = 2703.9524
$ws.Range("M7").Value = -2591.9524

$ws.Range("H20").Value = 1000000
$ws.Range("I20").Value = 1000000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1000000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -999774
$ws.Range("N20").ClearContents()

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H40").Value = 3602.697
$ws.Range("I40").Value = 2179.25
$ws.Range("K40").Value = 2179.25
$ws.Range("M40").Value = -2043.25

$ws.Range("H42").Value = 272006
$ws.Range("I42").Value = 352674.66
$ws.Range("J42").Value = 30000
$ws.Range("K42").Value = 352674.66
$ws.Range("L42").Value = 30000
$ws.Range("M42").Value = -352111.66
$ws.Range("N42").Value = -31126

$ws.Range("H46").Value = 1096.2941
$ws.Range("I46").Value = 651
$ws.Range("J46").Value = 2813.8572
$ws.Range("K46").Value = 651
$ws.Range("L46").Value = 2813.8572
$ws.Range("M46").Value = -463
$ws.Range("N46").Value = -3189.8572

$ws.Range("H49").Value = 272006
$ws.Range("I49").Value = 352674.66
$ws.Range("J49").Value = 30000
$ws.Range("K49").Value = 352674.66
$ws.Range("L49").Value = 30000
$ws.Range("M49").Value = -352527.66
$ws.Range("N49").Value = -30294

$ws.Range("H68").Value = 6074.8
$ws.Range("J68").Value = 7382.5
$ws.Range("L68").Value = 7382.5
$ws.Range("N68").Value = -8880.5

$ws.Range("H71").Value = 6074.8
$ws.Range("J71").Value = 7382.5
$ws.Range("L71").Value = 36912.5
$ws.Range("N71").Value = -44400.5

$ws.Range("H87").Value = 500012580
$ws.Range("J87").Value = 1000000000
$ws.Range("L87").Value = 1000000000
$ws.Range("N87").Value = -1000002246

$ws.Range("H90").Value = 500012580
$ws.Range("J90").Value = 1000000000
$ws.Range("L90").Value = 3000000000
$ws.Range("N90").Value = -3000011232

$ws.Range("H93").Value = 3448.375
$ws.Range("I93").Value = 1626.375
$ws.Range("K93").Value = 1626.375
$ws.Range("M93").Value = -378.375

$ws.Range("H122").Value = 2557.0908
$ws.Range("I122").Value = 2054.625
$ws.Range("J122").Value = 3897
$ws.Range("K122").Value = 6163.875
$ws.Range("L122").Value = 11691
$ws.Range("M122").Value = -3713.875
$ws.Range("N122").Value = -16591

$ws.Range("H126").Value = 3016.52
$ws.Range("I126").Value = 2703.9524
$ws.Range("K126").Value = 8111.8572
$ws.Range("M126").Value = -5641.8572

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 2318.8333
$ws.Range("I132").Value = 2181.8948
$ws.Range("J132").Value = 2839.2
$ws.Range("K132").Value = 6545.6844
$ws.Range("L132").Value = 8517.599999999999
$ws.Range("M132").Value = -4015.6844
$ws.Range("N132").Value = -13577.6

$ws.Range("H140").Value = 57713.453
$ws.Range("J140").Value = 57713.453
$ws.Range("L140").Value = 57713.453
$ws.Range("N140").Value = -68073.45300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 52500
$ws.Range("I29").Value = 50000
$ws.Range("J29").Value = 55000
$ws.Range("K29").Value = 50000
$ws.Range("L29").Value = 55000
$ws.Range("M29").Value = -49710
$ws.Range("N29").Value = -55580

$ws.Range("H38").Value = 2510.3333
$ws.Range("I38").Value = 2028
$ws.Range("J38").Value = 3475
$ws.Range("K38").Value = 2028
$ws.Range("L38").Value = 3475
$ws.Range("M38").Value = -1555
$ws.Range("N38").Value = -4421

$ws.Range("H62").Value = 12207183
$ws.Range("J62").Value = 14193.719
$ws.Range("L62").Value = 14193.719
$ws.Range("N62").Value = -15441.719

$ws.Range("H65").Value = 12207183
$ws.Range("J65").Value = 14193.719
$ws.Range("L65").Value = 70968.595
$ws.Range("N65").Value = -77208.595

$ws.Range("H74").Value = 42474.75
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 42474.75
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 42474.75
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -44346.75

$ws.Range("H77").Value = 42474.75
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 42474.75
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 127424.25
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -136784.25

$ws.Range("H107").Value = 21739476
$ws.Range("I107").Value = 322.92856
$ws.Range("K107").Value = 968.78568
$ws.Range("M107").Value = 951.21432

$ws.Range("H109").Value = 124699.8
$ws.Range("J109").Value = 124699.8
$ws.Range("L109").Value = 124699.8
$ws.Range("N109").Value = -127473.8

$ws.Range("H122").Value = 455576.53
$ws.Range("I122").Value = 974465.0600000001
$ws.Range("K122").Value = 2923395.18
$ws.Range("M122").Value = -2920945.18

$ws.Range("H126").Value = 11114149
$ws.Range("I126").Value = 2932
$ws.Range("K126").Value = 8796
$ws.Range("M126").Value = -6326

$ws.Range("H132").Value = 13343588
$ws.Range("I132").Value = 16678792
$ws.Range("J132").Value = 2771.2
$ws.Range("K132").Value = 50036376
$ws.Range("L132").Value = 8313.599999999999
$ws.Range("M132").Value = -50033846
$ws.Range("N132").Value = -13373.6

$ws.Range("H136").Value = 3942.1785
$ws.Range("I136").Value = 1271.1666
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 3813.4998
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = -1263.4998
$ws.Range("N136").Value = -31350
